$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-89 down to 16-90.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new data record.
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(15, 3).Value = 'Ñuble'
$ws.Cells.Item(15, 4).Value = 44473
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100112045
$ws.Cells.Item(15, 7).Value = 'Zapallo'
$ws.Cells.Item(15, 8).Value = 'Camote'
$ws.Cells.Item(15, 9).Value = '1a (guarda)'
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 400
$ws.Cells.Item(15, 12).Value = 450
$ws.Cells.Item(15, 13).Value = 425
$ws.Cells.Item(15, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(15, 15).Value = 'Región del Maule'
$ws.Cells.Item(15, 16).Value = 425
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(15, 18).Value = 'Hortaliza'
